$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") date value from 45172 to 45175 for all data rows (2-261)
$ws.Range("C2:C261").Value = 45175
